$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Support for Hot-Swapping multiple levels without a restart"
# feature's EARNED value (C14) from 0 to 0.11, matching the VALUE (B14).
$ws.Range("C14").Value = 0.11

# Color the earned value red, like the other completed feature cells in
# column C (this matches existing style used elsewhere, e.g. C2:C11).
$ws.Range("C14").Font.Color = 255

# Move the active selection to G19 (as left by the author after editing).
$ws.Range("G19").Select()
